$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,@(2, 2, 9)
    ,@(2, 3, 2)
    ,@(2, 4, 11)
    ,@(2, 6, 3)
    ,@(2, 7, 11)
    ,@(2, 10, 9)
    ,@(2, 11, 5)
    ,@(2, 12, 2)
    ,@(2, 13, 11)
    ,@(2, 14, 2)
    ,@(2, 15, 2)
    ,@(2, 16, 2)
    ,@(2, 17, 2)
    ,@(3, 2, 3)
    ,@(3, 3, 1)
    ,@(3, 4, 2)
    ,@(3, 5, 3)
    ,@(3, 6, 9)
    ,@(3, 7, 2)
    ,@(3, 8, 3)
    ,@(3, 10, 5)
    ,@(3, 12, 11)
    ,@(3, 13, 2)
    ,@(3, 14, 11)
    ,@(3, 15, 11)
    ,@(3, 16, 11)
    ,@(3, 17, 11)
    ,@(4, 2, 11)
    ,@(4, 3, 3)
    ,@(4, 4, 3)
    ,@(4, 5, 5)
    ,@(4, 6, 11)
    ,@(4, 7, 3)
    ,@(4, 8, 5)
    ,@(4, 10, 3)
    ,@(4, 11, 7)
    ,@(4, 12, 3)
    ,@(4, 13, 3)
    ,@(4, 14, 9)
    ,@(4, 15, 3)
    ,@(4, 16, 3)
    ,@(4, 17, 3)
    ,@(5, 2, 2)
    ,@(5, 3, 9)
    ,@(5, 4, 9)
    ,@(5, 5, 4)
    ,@(5, 6, 5)
    ,@(5, 7, 9)
    ,@(5, 8, 4)
    ,@(5, 10, 6)
    ,@(5, 11, 1)
    ,@(5, 12, 9)
    ,@(5, 13, 10)
    ,@(5, 14, 3)
    ,@(5, 15, 9)
    ,@(5, 16, 9)
    ,@(5, 17, 9)
    ,@(6, 2, 5)
    ,@(6, 3, 5)
    ,@(6, 4, 5)
    ,@(6, 5, 1)
    ,@(6, 6, 2)
    ,@(6, 7, 5)
    ,@(6, 8, 1)
    ,@(6, 10, 7)
    ,@(6, 11, 3)
    ,@(6, 12, 5)
    ,@(6, 13, 9)
    ,@(6, 14, 5)
    ,@(6, 15, 10)
    ,@(6, 16, 5)
    ,@(6, 17, 10)
    ,@(7, 2, 0)
    ,@(7, 3, 0)
    ,@(7, 4, 0)
    ,@(7, 5, 2)
    ,@(7, 6, 12)
    ,@(7, 7, 0)
    ,@(7, 8, 2)
    ,@(7, 10, 4)
    ,@(7, 11, 0)
    ,@(7, 12, 0)
    ,@(7, 14, 12)
    ,@(7, 15, 14)
    ,@(7, 16, 0)
    ,@(7, 17, 14)
    ,@(8, 2, 1)
    ,@(8, 3, 1)
    ,@(8, 4, 1)
    ,@(8, 5, 7)
    ,@(8, 6, 0)
    ,@(8, 7, 1)
    ,@(8, 8, 7)
    ,@(8, 10, 8)
    ,@(8, 11, 4)
    ,@(8, 12, 1)
    ,@(8, 13, 12)
    ,@(8, 14, 10)
    ,@(8, 15, 12)
    ,@(8, 16, 1)
    ,@(8, 17, 12)
    ,@(9, 2, 7)
    ,@(9, 3, 2)
    ,@(9, 4, 12)
    ,@(9, 5, 0)
    ,@(9, 6, 1)
    ,@(9, 7, 12)
    ,@(9, 8, 0)
    ,@(9, 10, 0)
    ,@(9, 11, 8)
    ,@(9, 12, 12)
    ,@(9, 13, 4)
    ,@(9, 14, 0)
    ,@(9, 15, 5)
    ,@(9, 16, 12)
    ,@(9, 17, 8)
    ,@(10, 2, 12)
    ,@(10, 3, 7)
    ,@(10, 4, 7)
    ,@(10, 5, 11)
    ,@(10, 6, 7)
    ,@(10, 7, 7)
    ,@(10, 8, 11)
    ,@(10, 10, 2)
    ,@(10, 11, 6)
    ,@(10, 12, 7)
    ,@(10, 13, 8)
    ,@(10, 14, 14)
    ,@(10, 15, 8)
    ,@(10, 16, 7)
    ,@(10, 17, 4)
    ,@(11, 2, 14)
    ,@(11, 3, 4)
    ,@(11, 4, 14)
    ,@(11, 5, 10)
    ,@(11, 6, 14)
    ,@(11, 7, 14)
    ,@(11, 8, 10)
    ,@(11, 11, 9)
    ,@(11, 12, 14)
    ,@(11, 13, 6)
    ,@(11, 14, 1)
    ,@(11, 15, 6)
    ,@(11, 16, 10)
    ,@(11, 17, 6)
    ,@(12, 2, 10)
    ,@(12, 3, 4)
    ,@(12, 4, 10)
    ,@(12, 5, 9)
    ,@(12, 6, 4)
    ,@(12, 7, 10)
    ,@(12, 8, 9)
    ,@(12, 12, 10)
    ,@(12, 13, 7)
    ,@(12, 14, 7)
    ,@(12, 15, 4)
    ,@(12, 16, 14)
    ,@(13, 2, 4)
    ,@(13, 3, 0)
    ,@(13, 4, 4)
    ,@(13, 5, 13)
    ,@(13, 6, 10)
    ,@(13, 7, 4)
    ,@(13, 8, 13)
    ,@(13, 12, 4)
    ,@(13, 13, 5)
    ,@(13, 15, 0)
    ,@(13, 16, 4)
    ,@(13, 17, 7)
    ,@(14, 2, 6)
    ,@(14, 3, 6)
    ,@(14, 4, 8)
    ,@(14, 5, 14)
    ,@(14, 6, 8)
    ,@(14, 7, 8)
    ,@(14, 8, 14)
    ,@(14, 12, 8)
    ,@(14, 13, 0)
    ,@(14, 14, 8)
    ,@(14, 15, 7)
    ,@(14, 16, 8)
    ,@(14, 17, 0)
    ,@(15, 2, 8)
    ,@(15, 3, 8)
    ,@(15, 4, 6)
    ,@(15, 5, 8)
    ,@(15, 6, 6)
    ,@(15, 7, 6)
    ,@(15, 8, 8)
    ,@(15, 12, 6)
    ,@(15, 14, 6)
    ,@(15, 15, 1)
    ,@(15, 16, 6)
    ,@(15, 17, 1)
    ,@(16, 2, 13)
    ,@(16, 3, 3)
    ,@(16, 4, 13)
    ,@(16, 5, 6)
    ,@(16, 6, 13)
    ,@(16, 7, 13)
    ,@(16, 8, 6)
    ,@(16, 12, 13)
    ,@(16, 13, 13)
    ,@(16, 14, 13)
    ,@(16, 15, 13)
    ,@(16, 16, 13)
    ,@(16, 17, 13)
)

foreach ($change in $changes) {
    $r = $change[0]
    $c = $change[1]
    $v = $change[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Host "Applied $($changes.Count) cell updates"
